# Sheaft product catalogue import template - fix categories list
# (renames singular category labels to plural, merges duplicate
# data-validation rules, and refreshes the remembered cell selections).

$wb  = $excel.ActiveWorkbook
$wsCatalogue = $wb.Worksheets.Item("Catalogue")
$wsListes    = $wb.Worksheets.Item("Listes")

# --- "Listes" sheet: rename category labels to their plural form ---
# (Tableau5 / liste_categorie, column G)
$wsListes.Range("G2").Value = "Boissons"
$wsListes.Range("G6").Value = "Poissons"
$wsListes.Range("G7").Value = "Viandes"

# --- "Catalogue" sheet: merge the "Agriculture Bio" (I) and
# "Disponible à la vente" (J) boolean-list validations into one rule
# covering both columns ---
$rngBool = $wsCatalogue.Range("I2:J1048576")
$rngBool.Validation.Delete()
$rngBool.Validation.Add(3, 1, 1, "liste_bool")
$rngBool.Validation.InputMessage = "Cliquez sur la flèche à droite de la case"
$rngBool.Validation.ShowInput = $true
$rngBool.Validation.ShowError = $true

# --- "Catalogue" sheet: merge the "Prix HT" (F) and "Quantité" (C)
# positive-decimal validations ---
$rngPrice = $wsCatalogue.Range("F2:F1048576")
$rngQty   = $wsCatalogue.Range("C2:C1048576")
$rngDecimal = $excel.Union($rngPrice, $rngQty)
$rngDecimal.Validation.Delete()
$rngDecimal.Validation.Add(2, 1, 5, "0")
$rngDecimal.Validation.ErrorTitle = "Erreur"
$rngDecimal.Validation.ErrorMessage = "Attention, il faut saisir un nombre (entier ou décimal) supérieur à 0. Cliquez sur ""Rééssayer"" pour modifier votre saisie."
$rngDecimal.Validation.ShowInput = $true
$rngDecimal.Validation.ShowError = $true

# --- Restore the remembered selections on each sheet (Listes first, so
# that "Catalogue" ends up as the active/selected tab, matching the
# original workbook) ---
$wsListes.Range("D10").Select() | Out-Null
$wsCatalogue.Range("H6").Select() | Out-Null
